# Vendor emails workbook update:
#  - remove the "NEUCO, INC. (BH)" vendor row entirely
#  - add a "Purchaser" column (C) and a "PurchaserEmail" column (D)
#  - populate every data row with the purchaser name "Lee Forsythe" and
#    a mailto hyperlink to leeforsythe@rogerssupply.com in column D
#  - adjust the view (selection / top-left cell) and column widths

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the NEUCO, INC. (BH) row (original row 4) -- all rows below shift up.
$ws.Rows(4).Delete()

# After the delete, 17 data rows remain in rows 2..18 (row 1 is the header).
$lastRow = 18

# 2. Re-point the existing Email hyperlink range (it still spans the old
#    B3:B19 range after the row delete) to the correct, shifted B3:B18 range.
#    Do this before adding any other links so relationship ids line up
#    (rId1 = B2, rId2 = B3:B18, rId3.. = the new PurchaserEmail links).
$ws.Hyperlinks.Item(2).Delete()
$ws.Hyperlinks.Add($ws.Range("B3:B18"), "mailto:leeforsythe@rogerssupply.com", "", "", "leeforsythe@rogerssupply.com")
$ws.Range("B3:B18").Style = "Hyperlink"

# 3. New header cells.
$ws.Range("C1").Value = "Purchaser"
$ws.Range("D1").Value = "PurchaserEmail"

# 4. Populate Purchaser / PurchaserEmail for every data row, and hyperlink each
#    PurchaserEmail cell individually (mirrors how column B's link was built).
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("C$r").Value = "Lee Forsythe"
    $ws.Range("D$r").Value = "leeforsythe@rogerssupply.com"
    $ws.Hyperlinks.Add($ws.Range("D$r"), "mailto:leeforsythe@rogerssupply.com")
    $ws.Range("D$r").Style = "Hyperlink"
}

# 5. Column widths for the two new columns.
$ws.Columns.Item(3).ColumnWidth = 17.5
$ws.Columns.Item(4).ColumnWidth = 20.666666666666668

# 6. View state: scrolled down a bit with A7 selected.
$ws.Range("A7").Select()
$excel.ActiveWindow.ScrollRow = 4
